# "finish dev of switch production building type"
#
# The "houses" sheet had a STR_limitBy column (E) that listed the building
# that limited/gated production (townHall, lumbermill, mill, stoneMason,
# foundry). That column is removed entirely: the old STR_preCondition (F)
# and STR_desc (G) columns shift left to become E and F, and the now-unused
# STR_limitBy strings are dropped from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column E (STR_limitBy) - remaining cells to the right (old F, G)
# shift left automatically, same as Excel's Delete with xlShiftToLeft.
$ws.Range("E1:E6").Delete()

# Leave the cursor where the author left it after finishing the edit.
$ws.Range("E5").Select()
